$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.331.69'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '2.648.80'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.91'
$ws.Range("E5").Value = '  -1.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.09'
$ws.Range("E6").Value = '  -2.91%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -0.80%  '
$ws.Range("D9").Value = '2.647.10'
$ws.Range("E9").Value = '  +0.17%  '
$ws.Range("E10").Value = '  -2.68%  '
$ws.Range("E11").Value = '  +1.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.357'
$ws.Range("E12").Value = '  -0.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.97'
$ws.Range("E13").Value = '  -1.80%  '
$ws.Range("D14").Value = '3.133.30'
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("E15").Value = '  -2.57%  '
$ws.Range("D16").Value = '72.203.69'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("E17").Value = '  -2.24%  '
$ws.Range("D18").Value = '2.658.56'
$ws.Range("E18").Value = '  +0.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.31'
$ws.Range("E19").Value = '  +2.46%  '
$ws.Range("E20").Value = '  +0.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '371.18'
$ws.Range("E21").Value = '  -2.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.16'
$ws.Range("E22").Value = '  -0.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.06'
$ws.Range("E23").Value = '  -0.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.48'
$ws.Range("E24").Value = '  -2.36%  '
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("E26").Value = '  -3.15%  '
$ws.Range("E27").Value = '  -4.05%  '
$ws.Range("D28").Value = '2.785.77'
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.29%  '
$ws.Range("D30").Value = '0.0₃0955'
$ws.Range("E30").Value = '  -0.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.07'
$ws.Range("E31").Value = '  -1.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '497.46'
$ws.Range("E32").Value = '  -5.13%  '
$ws.Range("E33").Value = '  -2.65%  '
$ws.Range("E34").Value = '  -1.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.95'
$ws.Range("E36").Value = '  -1.39%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.117'
$ws.Range("E37").Value = '  +3.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.38'
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("E39").Value = '  -0.98%  '
$ws.Range("E40").Value = '  -3.28%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("E42").Value = '  -6.86%  '
$ws.Range("E43").Value = '  -3.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.91'
$ws.Range("E44").Value = '  -3.62%  '
$ws.Range("E45").Value = '  -1.49%  '
$ws.Range("E46").Value = '  -0.55%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '153.97'
$ws.Range("E47").Value = '  +1.47%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.552'
$ws.Range("E48").Value = '  +0.93%  '
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.68'
$ws.Range("E49").Value = '  -1.20%  '
$ws.Range("E50").Value = '  -0.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0749'
$ws.Range("E51").Value = '  -1.19%  '
